# Insert a new data row at row 65 (shifting existing rows 65-184 down to 66-185)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(65).Insert()

$ws.Cells.Item(65, 1).Value = 5
$ws.Cells.Item(65, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(65, 3).Value = "Maule"
$ws.Cells.Item(65, 4).Value = 44495
$ws.Cells.Item(65, 5).Value = 7
$ws.Cells.Item(65, 6).Value = 100112006
$ws.Cells.Item(65, 7).Value = "Repollo"
$ws.Cells.Item(65, 8).Value = "Crespo record"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 5000
$ws.Cells.Item(65, 11).Value = 600
$ws.Cells.Item(65, 12).Value = 600
$ws.Cells.Item(65, 13).Value = 600
$ws.Cells.Item(65, 14).Value = "`$/unidad"
$ws.Cells.Item(65, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(65, 16).Value = 600
$ws.Cells.Item(65, 17).Value = 1
$ws.Cells.Item(65, 18).Value = "Hortaliza"
